$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 32

# Column A holds a date-like label ("01-07-2021") that must be stored as
# plain text (matching the existing rows), not auto-converted into an
# Excel date serial number. Enter it as a text formula first, then use
# Copy/PasteSpecial(values) to collapse it down to a literal text value
# without Excel re-inferring a date type or introducing a new cell style.
$cellA = $ws.Cells.Item($row, 1)
$cellA.Formula = "=""01-07-2021"""
$cellA.Copy()
$cellA.PasteSpecial(-4163)

$ws.Cells.Item($row, 2).Value = 152962
$ws.Cells.Item($row, 3).Value = 26115
$ws.Cells.Item($row, 4).Value = 16426
$ws.Cells.Item($row, 5).Value = 11548
$ws.Cells.Item($row, 6).Value = 9482
$ws.Cells.Item($row, 7).Value = 12288
$ws.Cells.Item($row, 8).Value = 27284
$ws.Cells.Item($row, 9).Value = 29699
$ws.Cells.Item($row, 10).Value = 20118
